# Auto-generated edit script applying numeric updates to the Leve profit
# tracking sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 437.9
$ws.Range("I38").Value = 264.33334
$ws.Range("J38").Value = 2000.0
$ws.Range("K38").Value = 793.0000200000001
$ws.Range("L38").Value = 6000.0
$ws.Range("M38").Value = -421.0000200000001
$ws.Range("N38").Value = -6744.0

$ws.Range("H74").Value = 3423.625
$ws.Range("I74").Value = 3055.5715
$ws.Range("K74").Value = 3055.5715
$ws.Range("M74").Value = -2119.5715

$ws.Range("H77").Value = 3423.625
$ws.Range("I77").Value = 3055.5715
$ws.Range("K77").Value = 15277.8575
$ws.Range("M77").Value = -10597.8575

$ws.Range("H86").Value = 95548.62
$ws.Range("I86").Value = 137650.22
$ws.Range("K86").Value = 137650.22
$ws.Range("M86").Value = -136527.22

$ws.Range("H89").Value = 95548.62
$ws.Range("I89").Value = 137650.22
$ws.Range("K89").Value = 688251.1
$ws.Range("M89").Value = -682635.1

$ws.Range("H98").Value = 1663.6666
$ws.Range("I98").Value = 1734.7142
$ws.Range("J98").Value = 1166.3334
$ws.Range("K98").Value = 1734.7142
$ws.Range("L98").Value = 1166.3334
$ws.Range("M98").Value = -236.7141999999999
$ws.Range("N98").Value = -4162.3334

$ws.Range("H122").Value = 1663.6666
$ws.Range("I122").Value = 1734.7142
$ws.Range("J122").Value = 1166.3334
$ws.Range("K122").Value = 5204.142599999999
$ws.Range("L122").Value = 3499.0002
$ws.Range("M122").Value = -2754.142599999999
$ws.Range("N122").Value = -8399.0002

$ws.Range("H127").Value = 1831.3334
$ws.Range("I127").Value = 1708.2307
$ws.Range("K127").Value = 5124.6921
$ws.Range("M127").Value = -164.6921000000002

$ws.Range("H131").Value = 1876.2106
$ws.Range("I131").Value = 788.0
$ws.Range("K131").Value = 2364.0
$ws.Range("M131").Value = 2676.0

$ws.Range("H132").Value = 601.8182
$ws.Range("I132").Value = 549.375
$ws.Range("K132").Value = 1648.125
$ws.Range("M132").Value = 881.875

$ws.Range("H137").Value = 2240.077
$ws.Range("I137").Value = 1488.6666
$ws.Range("J137").Value = 2465.5
$ws.Range("K137").Value = 4465.9998
$ws.Range("L137").Value = 7396.5
$ws.Range("M137").Value = -1915.9998
$ws.Range("N137").Value = -12496.5

$ws.Range("H138").Value = 2138.7932
$ws.Range("I138").Value = 1940.5667
$ws.Range("J138").Value = 2351.1785
$ws.Range("K138").Value = 5821.7001
$ws.Range("L138").Value = 7053.5355
$ws.Range("M138").Value = -681.7001
$ws.Range("N138").Value = -17333.5355

$ws.Range("H141").Value = 2591.1904
$ws.Range("I141").Value = 1028.6
$ws.Range("J141").Value = 6497.6665
$ws.Range("K141").Value = 3085.8
$ws.Range("L141").Value = 19492.9995
$ws.Range("M141").Value = 2094.2
$ws.Range("N141").Value = -29852.9995


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2241.1309
$ws.Range("I32").Value = 1805.5062
$ws.Range("J32").Value = 14003.0
$ws.Range("K32").Value = 1805.5062
$ws.Range("L32").Value = 14003.0
$ws.Range("M32").Value = -1518.5062
$ws.Range("N32").Value = -14577.0

$ws.Range("H102").Value = 2438.25
$ws.Range("I102").Value = 2438.25
$ws.Range("K102").Value = 2438.25
$ws.Range("M102").Value = -816.25

$ws.Range("H132").Value = 2249.4736
$ws.Range("I132").Value = 1488.7858
$ws.Range("J132").Value = 4379.4
$ws.Range("K132").Value = 4466.357400000001
$ws.Range("L132").Value = 13138.2
$ws.Range("M132").Value = -1936.357400000001
$ws.Range("N132").Value = -18198.2


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1746.0358
$ws.Range("I20").Value = 1616.5217
$ws.Range("J20").Value = 2341.8
$ws.Range("K20").Value = 1616.5217
$ws.Range("L20").Value = 2341.8
$ws.Range("M20").Value = -1369.5217
$ws.Range("N20").Value = -2835.8

$ws.Range("H86").Value = 170514.75
$ws.Range("I86").Value = 25000.0
$ws.Range("J86").Value = 183743.36
$ws.Range("K86").Value = 25000.0
$ws.Range("L86").Value = 183743.36
$ws.Range("M86").Value = -23877.0
$ws.Range("N86").Value = -185989.36

$ws.Range("H89").Value = 170514.75
$ws.Range("I89").Value = 25000.0
$ws.Range("J89").Value = 183743.36
$ws.Range("K89").Value = 125000.0
$ws.Range("L89").Value = 918716.7999999999
$ws.Range("M89").Value = -119384.0
$ws.Range("N89").Value = -929948.7999999999

$ws.Range("H105").Value = 6898481.0
$ws.Range("I105").Value = 2131.375
$ws.Range("J105").Value = 40000960.0
$ws.Range("K105").Value = 2131.375
$ws.Range("L105").Value = 40000960.0
$ws.Range("M105").Value = -384.375
$ws.Range("N105").Value = -40004454.0

$ws.Range("H107").Value = 2181.0
$ws.Range("I107").Value = 1728.2858
$ws.Range("J107").Value = 2973.25
$ws.Range("K107").Value = 1728.2858
$ws.Range("L107").Value = 2973.25
$ws.Range("M107").Value = 191.7141999999999
$ws.Range("N107").Value = -6813.25

$ws.Range("H134").Value = 6694.8667
$ws.Range("I134").Value = 7435.2915
$ws.Range("K134").Value = 22305.8745
$ws.Range("M134").Value = -19770.8745


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1943.5333
$ws.Range("I31").Value = 1781.625
$ws.Range("K31").Value = 1781.625
$ws.Range("M31").Value = -1486.625

$ws.Range("H34").Value = 1943.5333
$ws.Range("I34").Value = 1781.625
$ws.Range("K34").Value = 1781.625
$ws.Range("M34").Value = -1579.625

$ws.Range("H39").Value = 2000.0
$ws.Range("I39").Value = 2000.0
$ws.Range("K39").Value = 2000.0
$ws.Range("M39").Value = -1609.0

$ws.Range("H49").Value = 2000.0
$ws.Range("I49").Value = 2000.0
$ws.Range("K49").Value = 2000.0
$ws.Range("M49").Value = -1818.0

$ws.Range("H86").Value = 76924710.0
$ws.Range("I86").Value = 100001200.0
$ws.Range("K86").Value = 100001200.0
$ws.Range("M86").Value = -100000077.0

$ws.Range("H89").Value = 76924710.0
$ws.Range("I89").Value = 100001200.0
$ws.Range("K89").Value = 500006000.0
$ws.Range("M89").Value = -500000384.0

$ws.Range("H99").Value = 1884.55
$ws.Range("I99").Value = 1907.6666
$ws.Range("J99").Value = 1849.875
$ws.Range("K99").Value = 1907.6666
$ws.Range("L99").Value = 1849.875
$ws.Range("M99").Value = -409.6666
$ws.Range("N99").Value = -4845.875

$ws.Range("H105").Value = 1384.4375
$ws.Range("I105").Value = 1080.6428
$ws.Range("K105").Value = 1080.6428
$ws.Range("M105").Value = 666.3571999999999

$ws.Range("H107").Value = 916.375
$ws.Range("I107").Value = 918.8
$ws.Range("K107").Value = 918.8
$ws.Range("M107").Value = 1001.2

$ws.Range("H126").Value = 1884.55
$ws.Range("I126").Value = 1907.6666
$ws.Range("J126").Value = 1849.875
$ws.Range("K126").Value = 5722.9998
$ws.Range("L126").Value = 5549.625
$ws.Range("M126").Value = -3252.9998
$ws.Range("N126").Value = -10489.625

$ws.Range("H132").Value = 1996.3422
$ws.Range("I132").Value = 1139.409
$ws.Range("K132").Value = 3418.227
$ws.Range("M132").Value = -888.2270000000003


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 6183517.0
$ws.Range("J131").Value = 11205.806
$ws.Range("L131").Value = 33617.41800000001
$ws.Range("N131").Value = -43697.41800000001


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H106").Value = 0.0
$ws.Range("J106").Value = 0.0
$ws.Range("L106").Value = 0.0
$ws.Range("N106").ClearContents()

$ws.Range("H110").Value = 89998.5
$ws.Range("J110").Value = 89998.5
$ws.Range("L110").Value = 89998.5
$ws.Range("N110").Value = -98178.5

$ws.Range("H113").Value = 1133.1666
$ws.Range("J113").Value = 1133.1666
$ws.Range("L113").Value = 1133.1666
$ws.Range("N113").Value = -5473.1666

$ws.Range("H132").Value = 2266217.2
$ws.Range("I132").Value = 4810342.0
$ws.Range("J132").Value = 4773.0
$ws.Range("K132").Value = 14431026.0
$ws.Range("L132").Value = 14319.0
$ws.Range("M132").Value = -14428496.0
$ws.Range("N132").Value = -19379.0


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 8420.0
$ws.Range("J34").Value = 8420.0
$ws.Range("L34").Value = 8420.0
$ws.Range("N34").Value = -8764.0

$ws.Range("H40").Value = 3745.4666
$ws.Range("I40").Value = 2408.7778
$ws.Range("J40").Value = 5750.5
$ws.Range("K40").Value = 2408.7778
$ws.Range("L40").Value = 5750.5
$ws.Range("M40").Value = -2272.7778
$ws.Range("N40").Value = -6022.5

$ws.Range("H68").Value = 1829.6364
$ws.Range("J68").Value = 2666.6667
$ws.Range("L68").Value = 2666.6667
$ws.Range("N68").Value = -4164.6667

$ws.Range("H71").Value = 1829.6364
$ws.Range("J71").Value = 2666.6667
$ws.Range("L71").Value = 13333.3335
$ws.Range("N71").Value = -20821.3335

$ws.Range("H132").Value = 2093.25
$ws.Range("I132").Value = 1417.6
$ws.Range("J132").Value = 2575.8572
$ws.Range("K132").Value = 4252.799999999999
$ws.Range("L132").Value = 7727.571599999999
$ws.Range("M132").Value = -1722.799999999999
$ws.Range("N132").Value = -12787.5716


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 14944.0
$ws.Range("J104").Value = 14944.0
$ws.Range("L104").Value = 14944.0
$ws.Range("N104").Value = -21932.0

$ws.Range("H132").Value = 1643.3572
$ws.Range("J132").Value = 4999.0
$ws.Range("L132").Value = 14997.0
$ws.Range("N132").Value = -20057.0

$ws.Range("H136").Value = 11575487.0
$ws.Range("I136").Value = 14621054.0
$ws.Range("J136").Value = 2331.0
$ws.Range("K136").Value = 43863162.0
$ws.Range("L136").Value = 6993.0
$ws.Range("M136").Value = -43860612.0
$ws.Range("N136").Value = -12093.0

